# Auto-generated edit script: updates cached market-price / profit figures
# in the Famfrit_Profits workbook (scheduled-runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
    # Row 19: H19,I19,K19,M19
    $ws.Range("H19").Value = 6206.154
    $ws.Range("I19").Value = 1398.7142
    $ws.Range("K19").Value = 1398.7142
    $ws.Range("M19").Value = -1223.7142
    # Row 76: H76,I76,J76,K76,L76,M76,N76
    $ws.Range("H76").Value = 7977.85
    $ws.Range("I76").Value = 8643.3125
    $ws.Range("J76").Value = 5316
    $ws.Range("K76").Value = 8643.3125
    $ws.Range("L76").Value = 5316
    $ws.Range("M76").Value = -8328.3125
    $ws.Range("N76").Value = -5946
    # Row 79: H79,I79,J79,K79,L79,M79,N79
    $ws.Range("H79").Value = 7977.85
    $ws.Range("I79").Value = 8643.3125
    $ws.Range("J79").Value = 5316
    $ws.Range("K79").Value = 8643.3125
    $ws.Range("L79").Value = 5316
    $ws.Range("M79").Value = -7551.3125
    $ws.Range("N79").Value = -7500
    # Row 113: H113,I113,K113,M113
    $ws.Range("H113").Value = 7570
    $ws.Range("I113").Value = 5950.1665
    $ws.Range("K113").Value = 5950.1665
    $ws.Range("M113").Value = -2696.1665
    # Row 121: H121,J121,L121,N121
    $ws.Range("H121").Value = 1550.3334
    $ws.Range("J121").Value = 1550.3334
    $ws.Range("L121").Value = 4651.0002
    $ws.Range("N121").Value = -8145.0002
    # Row 125: H125,I125,J125,K125,L125,M125,N125
    $ws.Range("H125").Value = 3080.9285
    $ws.Range("I125").Value = 1955
    $ws.Range("J125").Value = 3706.4443
    $ws.Range("K125").Value = 17595
    $ws.Range("L125").Value = 33357.9987
    $ws.Range("M125").Value = -15135
    $ws.Range("N125").Value = -38277.9987
    # Row 132: H132,I132,J132,K132,L132,M132,N132
    $ws.Range("H132").Value = 1358.3091
    $ws.Range("I132").Value = 1042.4615
    $ws.Range("J132").Value = 6833
    $ws.Range("K132").Value = 3127.3845
    $ws.Range("L132").Value = 20499
    $ws.Range("M132").Value = -597.3844999999997
    $ws.Range("N132").Value = -25559
    # Row 137: H137,I137,K137,M137
    $ws.Range("H137").Value = 86961656
    $ws.Range("I137").Value = 58828164
    $ws.Range("K137").Value = 176484492
    $ws.Range("M137").Value = -176481942
    # Row 138: H138,I138,J138,K138,L138,M138,N138
    $ws.Range("H138").Value = 4071919.8
    $ws.Range("I138").Value = 3734.4707
    $ws.Range("J138").Value = 5135906.5
    $ws.Range("K138").Value = 11203.4121
    $ws.Range("L138").Value = 15407719.5
    $ws.Range("M138").Value = -6063.4121
    $ws.Range("N138").Value = -15417999.5

$ws = $wb.Worksheets.Item("ARM")
    # Row 32: H32,I32,J32,K32,L32,M32,N32
    $ws.Range("H32").Value = 18191768
    $ws.Range("I32").Value = 20415802
    $ws.Range("J32").Value = 28832
    $ws.Range("K32").Value = 20415802
    $ws.Range("L32").Value = 28832
    $ws.Range("M32").Value = -20415515
    $ws.Range("N32").Value = -29406
    # Row 45: H45,I45,K45,M45
    $ws.Range("H45").Value = 1994.125
    $ws.Range("I45").Value = 1358.5
    $ws.Range("K45").Value = 1358.5
    $ws.Range("M45").Value = -981.5
    # Row 74: H74,I74,J74,K74,L74,M74,N74
    $ws.Range("H74").Value = 125141770
    $ws.Range("I74").Value = 143018880
    $ws.Range("J74").Value = 1999
    $ws.Range("K74").Value = 143018880
    $ws.Range("L74").Value = 1999
    $ws.Range("M74").Value = -143018006
    $ws.Range("N74").Value = -3747
    # Row 77: H77,I77,J77,K77,L77,M77,N77
    $ws.Range("H77").Value = 125141770
    $ws.Range("I77").Value = 143018880
    $ws.Range("J77").Value = 1999
    $ws.Range("K77").Value = 715094400
    $ws.Range("L77").Value = 9995
    $ws.Range("M77").Value = -715090032
    $ws.Range("N77").Value = -18731
    # Row 132: H132,I132,J132,K132,L132,M132,N132
    $ws.Range("H132").Value = 26318814
    $ws.Range("I132").Value = 3030.4849
    $ws.Range("J132").Value = 200002980
    $ws.Range("K132").Value = 9091.4547
    $ws.Range("L132").Value = 600008940
    $ws.Range("M132").Value = -6561.4547
    $ws.Range("N132").Value = -600014000
    # Row 133: H133,I133,J133,K133,L133,M133,N133
    $ws.Range("H133").Value = 84093.414
    $ws.Range("I133").Value = 64137
    $ws.Range("J133").Value = 85907.63
    $ws.Range("K133").Value = 64137
    $ws.Range("L133").Value = 85907.63
    $ws.Range("M133").Value = -61607
    $ws.Range("N133").Value = -90967.63

$ws = $wb.Worksheets.Item("BSM")
    # Row 105: H105,I105,K105,M105
    $ws.Range("H105").Value = 12490.1
    $ws.Range("I105").Value = 12490.1
    $ws.Range("K105").Value = 12490.1
    $ws.Range("M105").Value = -10743.1

$ws = $wb.Worksheets.Item("CRP")
    # Row 22: H22,I22,J22,K22,L22,M22,N22
    $ws.Range("H22").Value = 4037.6072
    $ws.Range("I22").Value = 5613.579
    $ws.Range("J22").Value = 710.55554
    $ws.Range("K22").Value = 5613.579
    $ws.Range("L22").Value = 710.55554
    $ws.Range("M22").Value = -5263.579
    $ws.Range("N22").Value = -1410.55554
    # Row 31: H31,I31,J31,K31,L31,M31,N31
    $ws.Range("H31").Value = 24395066
    $ws.Range("I31").Value = 3633.276
    $ws.Range("J31").Value = 83341030
    $ws.Range("K31").Value = 3633.276
    $ws.Range("L31").Value = 83341030
    $ws.Range("M31").Value = -3338.276
    $ws.Range("N31").Value = -83341620
    # Row 34: H34,I34,J34,K34,L34,M34,N34
    $ws.Range("H34").Value = 24395066
    $ws.Range("I34").Value = 3633.276
    $ws.Range("J34").Value = 83341030
    $ws.Range("K34").Value = 3633.276
    $ws.Range("L34").Value = 83341030
    $ws.Range("M34").Value = -3431.276
    $ws.Range("N34").Value = -83341434
    # Row 58: H58,J58,L58,N58
    $ws.Range("H58").Value = 1866.4546
    $ws.Range("J58").Value = 1724.6666
    $ws.Range("L58").Value = 1724.6666
    $ws.Range("N58").Value = -2130.6666
    # Row 69: H69,I69,K69,M69
    $ws.Range("H69").Value = 135428.42
    $ws.Range("I69").Value = 95599.8
    $ws.Range("K69").Value = 95599.8
    $ws.Range("M69").Value = -94850.8
    # Row 72: H72,I72,K72,M72
    $ws.Range("H72").Value = 135428.42
    $ws.Range("I72").Value = 95599.8
    $ws.Range("K72").Value = 286799.4
    $ws.Range("M72").Value = -283055.4
    # Row 105: H105,I105,K105,M105
    $ws.Range("H105").Value = 18143.857
    $ws.Range("I105").Value = 4332.3335
    $ws.Range("K105").Value = 4332.3335
    $ws.Range("M105").Value = -2585.3335
    # Row 134: H134,I134,K134,M134
    $ws.Range("H134").Value = 2370.8
    $ws.Range("I134").Value = 2459
    $ws.Range("K134").Value = 7377
    $ws.Range("M134").Value = -4842
    # Row 136: H136,J136,L136,N136
    $ws.Range("H136").Value = 1866.4546
    $ws.Range("J136").Value = 1724.6666
    $ws.Range("L136").Value = 5173.9998
    $ws.Range("N136").Value = -10273.9998

$ws = $wb.Worksheets.Item("CUL")
    # Row 34: H34,J34,L34,N34
    $ws.Range("H34").Value = 2462
    $ws.Range("J34").Value = 3278
    $ws.Range("L34").Value = 9834
    $ws.Range("N34").Value = -10002
    # Row 55: H55,J55,L55,N55
    $ws.Range("H55").Value = 960.35
    $ws.Range("J55").Value = 1005.8823
    $ws.Range("L55").Value = 3017.6469
    $ws.Range("N55").Value = -3371.6469
    # Row 92: H92,I92,J92,K92,L92,M92,N92
    $ws.Range("H92").Value = 287.25
    $ws.Range("I92").Value = 252
    $ws.Range("J92").Value = 322.5
    $ws.Range("K92").Value = 756
    $ws.Range("L92").Value = 967.5
    $ws.Range("M92").Value = 492
    $ws.Range("N92").Value = -3463.5
    # Row 97: H97,I97,J97,K97,L97,M97,N97
    $ws.Range("H97").Value = 457.16666
    $ws.Range("I97").Value = 481.66666
    $ws.Range("J97").Value = 432.66666
    $ws.Range("K97").Value = 1444.99998
    $ws.Range("L97").Value = 1297.99998
    $ws.Range("M97").Value = -948.9999800000001
    $ws.Range("N97").Value = -2289.99998
    # Row 107: H107,J107,L107,N107
    $ws.Range("H107").Value = 1027
    $ws.Range("J107").Value = 1077.6428
    $ws.Range("L107").Value = 3232.9284
    $ws.Range("N107").Value = -7072.928400000001
    # Row 120: H120,I120,K120,M120
    $ws.Range("H120").Value = 4335.6665
    $ws.Range("I120").Value = 4335.6665
    $ws.Range("K120").Value = 13006.9995
    $ws.Range("M120").Value = -8168.999500000002
    # Row 121: H121,J121,L121,N121
    $ws.Range("H121").Value = 5983641
    $ws.Range("J121").Value = 25926684
    $ws.Range("L121").Value = 77780052
    $ws.Range("N121").Value = -77782672

$ws = $wb.Worksheets.Item("GSM")
    # Row 122: H122,I122,K122,M122
    $ws.Range("H122").Value = 2689.4
    $ws.Range("I122").Value = 2024.4166
    $ws.Range("K122").Value = 6073.2498
    $ws.Range("M122").Value = -3623.2498
    # Row 126: H126,I126,J126,K126,L126,M126,N126
    $ws.Range("H126").Value = 7699525
    $ws.Range("I126").Value = 4355294.5
    $ws.Range("J126").Value = 12506856
    $ws.Range("K126").Value = 13065883.5
    $ws.Range("L126").Value = 37520568
    $ws.Range("M126").Value = -13063413.5
    $ws.Range("N126").Value = -37525508
    # Row 132: H132,I132,J132,K132,L132,M132,N132
    $ws.Range("H132").Value = 3625.8408
    $ws.Range("I132").Value = 3252.5642
    $ws.Range("J132").Value = 6537.4
    $ws.Range("K132").Value = 9757.692599999998
    $ws.Range("L132").Value = 19612.2
    $ws.Range("M132").Value = -7227.692599999998
    $ws.Range("N132").Value = -24672.2

$ws = $wb.Worksheets.Item("LTW")
    # Row 22: H22,I22,J22,K22,L22,M22,N22
    $ws.Range("H22").Value = 3556.6667
    $ws.Range("I22").Value = 2286.75
    $ws.Range("J22").Value = 4338.154
    $ws.Range("K22").Value = 2286.75
    $ws.Range("L22").Value = 4338.154
    $ws.Range("M22").Value = -1991.75
    $ws.Range("N22").Value = -4928.154
    # Row 27: H27,I27,J27,K27,L27,M27,N27
    $ws.Range("H27").Value = 3556.6667
    $ws.Range("I27").Value = 2286.75
    $ws.Range("J27").Value = 4338.154
    $ws.Range("K27").Value = 2286.75
    $ws.Range("L27").Value = 4338.154
    $ws.Range("M27").Value = -2179.75
    $ws.Range("N27").Value = -4552.154
    # Row 63: H63,I63,J63,K63,L63,M63,N63
    $ws.Range("H63").Value = 56999
    $ws.Range("I63").Value = 57000
    $ws.Range("J63").Value = 56998
    $ws.Range("K63").Value = 57000
    $ws.Range("L63").Value = 56998
    $ws.Range("M63").Value = -56251
    $ws.Range("N63").Value = -58496
    # Row 66: H66,I66,J66,K66,L66,M66,N66
    $ws.Range("H66").Value = 56999
    $ws.Range("I66").Value = 57000
    $ws.Range("J66").Value = 56998
    $ws.Range("K66").Value = 171000
    $ws.Range("L66").Value = 170994
    $ws.Range("M66").Value = -167256
    $ws.Range("N66").Value = -178482
    # Row 93: H93,I93,J93,K93,L93,M93,N93
    $ws.Range("H93").Value = 2753.6667
    $ws.Range("I93").Value = 1755.5
    $ws.Range("J93").Value = 4750
    $ws.Range("K93").Value = 1755.5
    $ws.Range("L93").Value = 4750
    $ws.Range("M93").Value = -507.5
    $ws.Range("N93").Value = -7246
    # Row 122: H122,I122,J122,K122,L122,M122,N122
    $ws.Range("H122").Value = 7099.5713
    $ws.Range("I122").Value = 4999.5
    $ws.Range("J122").Value = 7939.6
    $ws.Range("K122").Value = 14998.5
    $ws.Range("L122").Value = 23818.8
    $ws.Range("M122").Value = -12548.5
    $ws.Range("N122").Value = -28718.8
    # Row 131: H131,J131,L131,N131
    $ws.Range("H131").Value = 53763.5
    $ws.Range("J131").Value = 53763.5
    $ws.Range("L131").Value = 53763.5
    $ws.Range("N131").Value = -63843.5
    # Row 136: H136,I136,J136,K136,L136,M136,N136
    $ws.Range("H136").Value = 8735.474
    $ws.Range("I136").Value = 8266.733
    $ws.Range("J136").Value = 10493.25
    $ws.Range("K136").Value = 24800.199
    $ws.Range("L136").Value = 31479.75
    $ws.Range("M136").Value = -22250.199
    $ws.Range("N136").Value = -36579.75

$ws = $wb.Worksheets.Item("WVR")
    # Row 59: H59,J59,L59,N59
    $ws.Range("H59").Value = 40000
    $ws.Range("J59").Value = 40000
    $ws.Range("L59").Value = 40000
    $ws.Range("N59").Value = -41476
    # Row 136: H136,I136,K136,M136
    $ws.Range("H136").Value = 1502.9
    $ws.Range("I136").Value = 1492.1111
    $ws.Range("K136").Value = 4476.3333
    $ws.Range("M136").Value = -1926.3333

